# Remove duplicate / redundant OD-pair rows from the straight-line distance
# table (rows 74-100 originally). Row 74 duplicates row 73 (reverse pair),
# and among rows 75-100 every reverse-duplicate pair beyond the first
# occurrence is dropped. Deleting bottom-to-top keeps row numbers stable
# for rows not yet processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(100, 98, 96, 94, 92, 90, 88, 86, 85, 83, 80, 78, 76, 74)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
